$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Header text updates ---
$ws.Range("A8").Value = "Volume 32   Number  26"
$ws.Range("C9").Value = "Report Covering the Week  6/23/2025  Through  6/29/2025"

# --- Cells whose type/style changes (text <-> number): copy from a matching template cell ---
$ws.Range("C14").Copy($ws.Range("D14"))
$ws.Range("L14").Copy($ws.Range("E14"))
$ws.Range("F14").Copy($ws.Range("C15"))
$ws.Range("I14").Copy($ws.Range("D20"))
$ws.Range("H15").Copy($ws.Range("E20"))
$ws.Range("G14").Copy($ws.Range("C22"))
$ws.Range("F23").Copy($ws.Range("D22"))
$ws.Range("M29").Copy($ws.Range("E22"))
$ws.Range("D15").Copy($ws.Range("C27"))
$ws.Range("J23").Copy($ws.Range("D33"))
$ws.Range("H14").Copy($ws.Range("E33"))

# --- Plain numeric value updates (style unchanged) ---
$ws.Range("F15").Value = 3
$ws.Range("C16").Value = 9
$ws.Range("E16").Value = 200
$ws.Range("F16").Value = 22
$ws.Range("G16").Value = 12
$ws.Range("H16").Value = 83.333333333333
$ws.Range("I16").Value = 122
$ws.Range("J16").Value = 95
$ws.Range("K16").Value = 28.421052631578
$ws.Range("L16").Value = 130.188679245283
$ws.Range("M16").Value = 144
$ws.Range("N16").Value = -80.996884735202
$ws.Range("D17").Value = 10
$ws.Range("E17").Value = -40
$ws.Range("F17").Value = 26
$ws.Range("G17").Value = 32
$ws.Range("H17").Value = -18.75
$ws.Range("I17").Value = 168
$ws.Range("J17").Value = 142
$ws.Range("K17").Value = 18.309859154929
$ws.Range("L17").Value = 71.428571428571
$ws.Range("M17").Value = 158.461538461538
$ws.Range("N17").Value = -31.147540983606
$ws.Range("C18").Value = 5
$ws.Range("D18").Value = 1
$ws.Range("E18").Value = 400
$ws.Range("F18").Value = 14
$ws.Range("G18").Value = 12
$ws.Range("H18").Value = 16.666666666666
$ws.Range("I18").Value = 101
$ws.Range("J18").Value = 96
$ws.Range("K18").Value = 5.208333333333
$ws.Range("L18").Value = 40.277777777777
$ws.Range("M18").Value = 12.222222222222
$ws.Range("N18").Value = -90.560747663551
$ws.Range("C19").Value = 32
$ws.Range("D19").Value = 44
$ws.Range("E19").Value = -27.272727272727
$ws.Range("F19").Value = 112
$ws.Range("G19").Value = 140
$ws.Range("H19").Value = -20
$ws.Range("I19").Value = 841
$ws.Range("J19").Value = 880
$ws.Range("K19").Value = -4.431818181818
$ws.Range("L19").Value = -7.276736493936
$ws.Range("M19").Value = 10.367454068241
$ws.Range("N19").Value = -76.073968705547
$ws.Range("C20").Value = 2
$ws.Range("G20").Value = 5
$ws.Range("H20").Value = -20
$ws.Range("I20").Value = 32
$ws.Range("K20").Value = 45.454545454545
$ws.Range("L20").Value = -8.571428571428
$ws.Range("M20").Value = 45.454545454545
$ws.Range("N20").Value = -87.044534412955
$ws.Range("C21").Value = 54
$ws.Range("D21").Value = 58
$ws.Range("E21").Value = -6.896551724137
$ws.Range("F21").Value = 181
$ws.Range("G21").Value = 202
$ws.Range("H21").Value = -10.396039603960
$ws.Range("I21").Value = 1285
$ws.Range("J21").Value = 1241
$ws.Range("K21").Value = 3.545527800161
$ws.Range("L21").Value = 9.829059829059
$ws.Range("M21").Value = 28.371628371628
$ws.Range("N21").Value = -77.621037965865
$ws.Range("F22").Value = 3
$ws.Range("H22").Value = -25
$ws.Range("I22").Value = 37
$ws.Range("J22").Value = 46
$ws.Range("K22").Value = -19.565217391304
$ws.Range("L22").Value = -7.5
$ws.Range("M22").Value = 42.307692307692
$ws.Range("C24").Value = 55
$ws.Range("D24").Value = 78
$ws.Range("E24").Value = -29.487179487179
$ws.Range("F24").Value = 189
$ws.Range("G24").Value = 270
$ws.Range("H24").Value = -30
$ws.Range("I24").Value = 1252
$ws.Range("J24").Value = 1469
$ws.Range("K24").Value = -14.771953710006
$ws.Range("L24").Value = -4.061302681992
$ws.Range("M24").Value = 42.11123723042
$ws.Range("C25").Value = 39
$ws.Range("D25").Value = 83
$ws.Range("E25").Value = -53.012048192771
$ws.Range("F25").Value = 155
$ws.Range("G25").Value = 246
$ws.Range("H25").Value = -36.991869918699
$ws.Range("I25").Value = 1154
$ws.Range("J25").Value = 1424
$ws.Range("K25").Value = -18.960674157303
$ws.Range("L25").Value = -12.839879154078
$ws.Range("C26").Value = 11
$ws.Range("D26").Value = 17
$ws.Range("E26").Value = -35.294117647058
$ws.Range("F26").Value = 59
$ws.Range("G26").Value = 50
$ws.Range("H26").Value = 18
$ws.Range("I26").Value = 355
$ws.Range("J26").Value = 335
$ws.Range("K26").Value = 5.970149253731
$ws.Range("L26").Value = 2.601156069364
$ws.Range("M26").Value = 51.063829787234
$ws.Range("F27").Value = 3
$ws.Range("H27").Value = 200
$ws.Range("C28").Value = 3
$ws.Range("D28").Value = 1
$ws.Range("E28").Value = 200
$ws.Range("F28").Value = 9
$ws.Range("G28").Value = 11
$ws.Range("H28").Value = -18.181818181818
$ws.Range("I28").Value = 63
$ws.Range("J28").Value = 42
$ws.Range("K28").Value = 50
$ws.Range("L28").Value = 57.5
$ws.Range("L31").Value = -23.076923076923
$ws.Range("G33").Value = 1
$ws.Range("J33").Value = 2
$ws.Range("K33").Value = 50
